$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 3 (formats/styles are inherited from the row
# that shifts down, matching the original row 3 which becomes row 4).
$ws.Rows("3:3").Insert()

# Populate the new row's Model / Production Performance columns.
$ws.Range("G3").Value = "gpt-4o-mini"
$ws.Range("H3").Value = "16.8s"

# Resize Table1 to include the newly inserted row.
[void]$ws.ListObjects.Item(1).Resize($ws.Range("A1:H10"))

# Update the active selection to match the saved state.
[void]$ws.Range("I4").Select()
